$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sort the product table (rows 2:24) alphabetically by brand (col B), ---
# --- same as the original "Shogun/Ghost" catalogue being re-grouped      ---
$dataRange = $ws.Range("A2:G24")
$sortKey = $ws.Range("B2:B24")
$dataRange.Sort($sortKey)

# --- Targeted cell edits on top of the new (sorted) layout ---
# Shorten / rename a few product + supplier names
$ws.Cells.Item(2,5).Value = "ghost importados"
$ws.Cells.Item(3,5).Value = "ghost importados"

$ws.Cells.Item(4,5).Value = "SHOGUN IMPORTS"
$ws.Cells.Item(5,5).Value = "SHOGUN IMPORTS"

$ws.Cells.Item(9,5).Value = "Ghost Importados"
$ws.Cells.Item(10,5).Value = "Ghost Importados"
$ws.Cells.Item(11,5).Value = "Ghost Importados"
$ws.Cells.Item(12,5).Value = "Ghost Importados"
$ws.Cells.Item(13,5).Value = "Ghost Importados"

$ws.Cells.Item(14,1).Value = "CONSOLE MEGA DRIVE"
$ws.Cells.Item(14,5).Value = "Ghost Importados"

$ws.Cells.Item(15,1).Value = "CONSOLE SATURN"
$ws.Cells.Item(15,5).Value = "Ghost Importados"

$ws.Cells.Item(16,1).Value = "CONSOLE DREAMCAST"
$ws.Cells.Item(16,5).Value = "Ghost Importados"

# --- Re-point the hyperlinks on column G to the rows they now belong to ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G17"), "mailto:shogun@imports.com")
$ws.Hyperlinks.Add($ws.Range("G9"), "mailto:ghost@importados.com.br")
$ws.Hyperlinks.Add($ws.Range("G3:G14"), "mailto:shogun@imports.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "shogun@imports.com")
$ws.Hyperlinks.Add($ws.Range("G16:G24"), "mailto:ghost@importados.com.br", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ghost@importados.com.br")

# --- Header row formatting: bold text on a grey fill ---
$headerRng = $ws.Range("A1:G1")
$headerRng.Font.Bold = $true
$headerRng.Interior.Color = 0xAAAAAE

# --- Column A width tweak ---
$ws.Columns.Item(1).ColumnWidth = 29

# --- Rename the main sheet and add the two extra (empty) sheets ---
$ws.Name = "consoles"
$s2 = $wb.Worksheets.Add($null, $ws)
$s2.Name = "Planilha1"
$s3 = $wb.Worksheets.Add($null, $s2)
$s3.Name = "Planilha2"

# --- Keep "consoles" as the active/selected tab, with the new selection ---
$ws.Activate()
$ws.Range("C30").Select()

Write-Host "edit complete"
